$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -0.0349
$ws.Range("C6").Value = -0.2623
$ws.Range("C14").Value = 0.0982
$ws.Range("C15").Value = 0.4627
$ws.Range("C16").Value = 0.7745
$ws.Range("C17").Value = 0.5919
$ws.Range("C19").Value = -0.0626
$ws.Range("C20").Value = 0.6402
